# Update the date placeholder on the notes master (datetimeFigureOut field)
$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
foreach ($sh in $nm.Shapes) {
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -eq "7/6/20") {
            $tr.Text = "10/1/20"
        }
    }
}

# Remove the "Dash-board" shape and expand the "Vault" shape into the freed space
$s = $p.Slides.Item(1)

$dashboard = $s.Shapes.Item("Rounded Rectangle 21")
$dashboard.Delete()

# PowerPoint's Shape.Left/Top/Width/Height are expressed in points (1 pt = 12700 EMU)
$vault = $s.Shapes.Item("Rounded Rectangle 33")
$vault.Left = 9482106 / 12700
$vault.Top = 7762234 / 12700
$vault.Width = 4129266 / 12700
$vault.Height = 1424957 / 12700
